$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the FilesTab Cypher query text in B4 ---
# Removed the "File Type" and "Breed" RETURN lines, and reflowed the
# indentation of the two lines that now follow the removed ones
# (Association / Diagnosis), matching the corrected ICDC script.
$newQuery = @'
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
WHERE demo.breed IN ['Samoyed']
OPTIONAL MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
WITH DISTINCT f, parent, c, demo, diag, s
RETURN  coalesce(f.file_name, '') AS `File Name`,
           coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(f.file_format, '') AS `Format`,
        coalesce(f.file_size, '') AS `Size`,
        coalesce(c.case_id, '') AS `Case ID`,
         coalesce(diag.disease_term,'') AS Diagnosis , 
        coalesce(s.clinical_study_designation,'') AS `Study Code`
'@
$ws.Range("B4").Value = $newQuery

# --- Row 4 is shorter now (2 fewer wrapped lines) so Excel re-fit the row height ---
$ws.Rows.Item(4).RowHeight = 217.5

# --- The sheet view had scrolled / the selection moved to B4 ---
$ws.Application.ActiveWindow.ScrollRow = 16
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("B4").Select()
